# Update "想去人数" (want-to-go count) figures for two exhibitions.
# These values are duplicated across the "展览" and "全部类型" worksheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F5").Value = 858
    $ws.Range("F7").Value = 424
}
